# Add season-record columns ("Wins", "Losses", "Ties") as AD, AE, AF
# to the roster sheet, matching the header style already used by the
# other header cells (e.g. AC1) and filling every data row (2-58) with
# the team's season record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1, bold /
# bordered / centered) onto the three new header cells before writing
# their text, so they pick up the same style as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# Season record is identical for every player on the roster.
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 73
    $ws.Cells.Item($r, 32).Value = 0
}
